$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.907.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.30"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4698"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3664"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07161"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9265"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.59"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07704"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.834.42"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.281"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.408"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.30"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008646"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.947.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.44"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.025"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.61"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.933"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.25"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.016"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.44"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.881"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08860"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.211"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.178"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7461"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.785"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.477"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.088"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01940"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.967"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05204"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5217"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.948"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1519"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.155"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.47"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4702"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.006"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.76"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.600"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.86"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06042"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8966"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.82%  "
